$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Notes" (sheet1)
# ---------------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item(1)
$wsNotes.Range("A2").Value = "This is an example of input data that should fail tests"
$wsNotes.Range("A3").Value = "Specific issue: counts table contains variant counts that sum to more than the total_num"
$wsNotes.Range("A4").Select()

# ---------------------------------------------------------------------------
# Sheet "studies" (sheet2)
# ---------------------------------------------------------------------------
$wsStudies = $wb.Worksheets.Item(2)
$wsStudies.Range("A1").Value = "study_id"
$wsStudies.Range("A2").Value = "study01"
$wsStudies.Range("B2").Value = "example name"
$wsStudies.Range("D2").Value = "Blaggs_etal"
$wsStudies.Range("E2").Value = 2024

# ---------------------------------------------------------------------------
# Sheet "surveys" (sheet3)
# ---------------------------------------------------------------------------
$wsSurveys = $wb.Worksheets.Item(3)

# Remove the extra example rows (rows 3-8), keep just one data row
$wsSurveys.Rows("3:8").Delete()

# Update header row text
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"

# Update data row (row 2)
$wsSurveys.Range("A2").Value = "study01"
$wsSurveys.Range("B2").Value = "S01"
$wsSurveys.Range("C2").Value = "Gambia"
$wsSurveys.Range("D2").Value = "example site"
$wsSurveys.Range("E2").Value = 0
$wsSurveys.Range("F2").Value = 0
$wsSurveys.Range("G2").Value = "example data"
$wsSurveys.Range("H2").Value = "2020-01-01"
$wsSurveys.Range("I2").Value = "2020-01-01"
$wsSurveys.Range("J2").Value = "2020-01-01"
$wsSurveys.Range("K2").Value = "example data"

# New header formatting: font size 12, black color for the whole header row
$wsSurveys.Range("A1:K1").Font.Size = 12
$wsSurveys.Range("A1:K1").Font.Color = 0
# H1:J1 keep a text number format (matches existing H:J column style)
$wsSurveys.Range("H1:J1").NumberFormat = "@"

$wsSurveys.Range("A1:K1").Select()

# ---------------------------------------------------------------------------
# Sheet "counts" (sheet4)
# ---------------------------------------------------------------------------
$wsCounts = $wb.Worksheets.Item(4)

# Remove the extra example rows (rows 4-8), keep just two data rows
$wsCounts.Rows("4:8").Delete()

$wsCounts.Range("A2").Value = "S01"
$wsCounts.Range("B2").Value = "crt:1_2_3:AAA;mdr1:1_2_3:AAA"
$wsCounts.Range("C2").Value = 5

$wsCounts.Range("A3").Value = "S01"
$wsCounts.Range("B3").Value = "crt:1_2_3:AAA;mdr1:1_2_3:AAC"
$wsCounts.Range("C3").Value = 6

$wsCounts.Columns("B:B").ColumnWidth = 32.998697916666664

$wsCounts.Range("D4").Select()

# ---------------------------------------------------------------------------
# Activate "studies" last so it becomes the workbook's active tab, matching
# the target activeTab="1" / tabSelected on sheet2.
# ---------------------------------------------------------------------------
$wsStudies.Activate()
$wsStudies.Range("A2").Select()
